$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Level" column (C) from "easy" to "moderate" for the rows below.
$rows = @(3, 4, 5, 7, 8, 9, 11)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "moderate"
}
